$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-13 (columns A-K)
$data = @(
    @(2240,  "18PG",  0, 9.3028, 1, 5,  20,  70, 160, 0.001327044311833706,  2.011390021328469),
    @(3690,  "18PG",  0, 9.3028, 1, 6,  150, 35, 40,  0.001326794255551539,  2.011387875095399),
    @(7653,  "18PG",  0, 9.3028, 1, 9,  200, 65, 5,   0.0006893382045958312, 2.005916589810046),
    @(747,   "18PG",  0, 9.3028, 1, 5,  70,  75, 7,   0.00006117747643906001, 2.000525086280276),
    @(10824, "18PG",  0, 9.3028, 1, 10, 40,  75, 30,  0.0003723426932675778, 2.003195817336316),
    @(77583, "DDAB",  1, 12.73,  0, 7,  200, 65, 9,   0.7638837941433332,    8.556414605132229),
    @(33757, "DOTAP", 1, 12.515, 2, 7,  50,  75, 3,   0.8007117726331708,    8.872509144510504),
    @(29038, "DOTAP", 1, 12.515, 2, 5,  40,  50, 4,   0.8348125341793452,    9.16519598086132),
    @(31284, "DOTAP", 1, 12.515, 2, 6,  150, 35, 3,   0.9024208400522759,    9.745478070168682),
    @(31499, "DOTAP", 1, 12.515, 2, 6,  450, 75, 4,   0.8339674949742872,    9.157943009364308),
    @(29214, "DOTAP", 1, 12.515, 2, 5,  250, 70, 5,   0.9008149036590295,    9.73169431810545),
    @(52416, "DSPC",  1, 6.12,   0, 9,  20,  25, 160, 0.001119367752479003,  2.009607533419527)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
    $ws.Cells.Item($row, 11).Value = $vals[10]
}
